# Updates crypto price/volume figures per the GitHub Actions refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.088.04'
$ws.Range("E2").Value = '  -2.31%  '
$ws.Range("D3").Value = '3.480.91'
$ws.Range("E3").Value = '  -1.19%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = "'" + '588.64'
$ws.Range("E5").Value = '  -3.15%  '
$ws.Range("D6").Value = "'" + '138.00'
$ws.Range("E6").Value = '  -3.78%  '
$ws.Range("D7").Value = '3.477.17'
$ws.Range("E7").Value = '  -1.26%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E9").Value = '  -4.49%  '
$ws.Range("E10").Value = '  -6.00%  '
$ws.Range("D11").Value = "'" + '7.22'
$ws.Range("E11").Value = '  -6.80%  '
$ws.Range("D12").Value = "'" + '0.381'
$ws.Range("E12").Value = '  -6.80%  '
$ws.Range("D13").Value = '4.065.47'
$ws.Range("E13").Value = '  -1.31%  '
$ws.Range("D14").Value = "'" + '0.0000182'
$ws.Range("E14").Value = '  -6.66%  '
$ws.Range("D15").Value = "'" + '26.62'
$ws.Range("E15").Value = '  -7.24%  '
$ws.Range("D16").Value = '3.452.61'
$ws.Range("E16").Value = '  -2.02%  '
$ws.Range("E17").Value = '  -1.35%  '
$ws.Range("D18").Value = '64.974.81'
$ws.Range("E18").Value = '  -2.26%  '
$ws.Range("D19").Value = "'" + '9.70'
$ws.Range("E19").Value = '  -9.93%  '
$ws.Range("D20").Value = "'" + '5.77'
$ws.Range("E20").Value = '  -6.37%  '
$ws.Range("D21").Value = "'" + '13.86'
$ws.Range("E21").Value = '  -5.65%  '
$ws.Range("D22").Value = "'" + '388.15'
$ws.Range("E22").Value = '  -8.36%  '
$ws.Range("E23").Value = '  -5.92%  '
$ws.Range("D24").Value = "'" + '1.00'
$ws.Range("E24").Value = '  -0.14%  '
$ws.Range("D25").Value = "'" + '72.45'
$ws.Range("E25").Value = '  -5.80%  '
$ws.Range("D26").Value = '3.621.42'
$ws.Range("E26").Value = '  -1.13%  '
$ws.Range("D27").Value = "'" + '5.75'
$ws.Range("E27").Value = '  -0.15%  '
$ws.Range("D28").Value = "'" + '0.0000109'
$ws.Range("E28").Value = '  -4.53%  '
$ws.Range("D29").Value = "'" + '0.999'
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("D30").Value = "'" + '7.41'
$ws.Range("E30").Value = '  -6.66%  '
$ws.Range("D31").Value = "'" + '8.21'
$ws.Range("E31").Value = '  -8.07%  '
$ws.Range("D32").Value = "'" + '2.22'
$ws.Range("E32").Value = '  -10.23%  '
$ws.Range("D33").Value = '3.493.79'
$ws.Range("E33").Value = '  -0.98%  '
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("E35").Value = '  -7.61%  '
$ws.Range("D36").Value = "'" + '23.00'
$ws.Range("E36").Value = '  -5.06%  '
$ws.Range("D37").Value = "'" + '171.87'
$ws.Range("E37").Value = '  -1.03%  '
$ws.Range("D38").Value = "'" + '1.21'
$ws.Range("E38").Value = '  -9.44%  '
$ws.Range("D39").Value = "'" + '6.89'
$ws.Range("E39").Value = '  -9.15%  '
$ws.Range("D40").Value = "'" + '1.47'
$ws.Range("E40").Value = '  -9.71%  '
$ws.Range("D41").Value = "'" + '4.76'
$ws.Range("E41").Value = '  -8.87%  '
$ws.Range("D42").Value = "'" + '0.0775'
$ws.Range("E42").Value = '  -4.87%  '
$ws.Range("D43").Value = "'" + '0.814'
$ws.Range("E43").Value = '  -4.70%  '
$ws.Range("D44").Value = "'" + '0.999'
$ws.Range("E44").Value = '  -0.08%  '
$ws.Range("D45").Value = "'" + '42.32'
$ws.Range("E45").Value = '  -6.98%  '
$ws.Range("D46").Value = "'" + '4.36'
$ws.Range("E46").Value = '  -12.73%  '
$ws.Range("D47").Value = "'" + '24.06'
$ws.Range("E47").Value = '  +5.35%  '
$ws.Range("D48").Value = "'" + '1.62'
$ws.Range("E48").Value = '  -8.84%  '
$ws.Range("D49").Value = "'" + '1.14'
$ws.Range("E49").Value = '  +1.04%  '
$ws.Range("D50").Value = "'" + '6.66'
$ws.Range("E50").Value = '  -5.90%  '
$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").Value = "'" + '2.07'
$ws.Range("E51").Value = '  -12.59%  '
